$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2, column A previously held a shared-string value (index 31); it was
# cleared out in the edit, leaving the cell blank.
$ws.Range("A2").ClearContents()

# Duplicate row 5 (the last data row) into four new rows, 7 through 10.
$srcRow = $ws.Range("A5:AG5")
$srcRow.Copy($ws.Range("A7:AG7"))
$srcRow.Copy($ws.Range("A8:AG8"))
$srcRow.Copy($ws.Range("A9:AG9"))
$srcRow.Copy($ws.Range("A10:AG10"))

# Leave the cursor where the author left it when they saved the file.
[void]$ws.Range("C17").Select()
